# The workbook originally has two sheets: "TMA map" (active/selected) and
# "er". This commit ("fix tab names in input TMAs") renames the second
# sheet from "er" to "ER", and the author was last looking at that sheet
# when they saved, so it becomes the selected/active tab instead of
# "TMA map".

$wb = $excel.ActiveWorkbook

# Rename the "er" sheet to "ER".
$erSheet = $wb.Worksheets.Item("er")
$erSheet.Name = "ER"

# Make "ER" the active sheet - this moves the active tab / selected tab
# from "TMA map" to "ER", matching the saved view state in the diff.
$erSheet.Activate()
